# Auto-generated edit script applying scheduled market-data refresh
# (H/I/J/K/L/M/N columns: currentAveragePrice*, LevePrice*, LeveProfit*)
$wb = $excel.ActiveWorkbook

$setOps = @(
    @{Sheet="ALC"; Cell="H33"; Value=728.2222},
    @{Sheet="ALC"; Cell="I33"; Value=110.1875},
    @{Sheet="ALC"; Cell="K33"; Value=110.1875},
    @{Sheet="ALC"; Cell="M33"; Value=118.8125},
    @{Sheet="ALC"; Cell="H43"; Value=4906.5},
    @{Sheet="ALC"; Cell="J43"; Value=5321.7144},
    @{Sheet="ALC"; Cell="L43"; Value=5321.7144},
    @{Sheet="ALC"; Cell="N43"; Value=-5459.7144},
    @{Sheet="ALC"; Cell="H55"; Value=620.53845},
    @{Sheet="ALC"; Cell="I55"; Value=338.6},
    @{Sheet="ALC"; Cell="J55"; Value=796.75},
    @{Sheet="ALC"; Cell="K55"; Value=338.6},
    @{Sheet="ALC"; Cell="L55"; Value=796.75},
    @{Sheet="ALC"; Cell="M55"; Value=-124.6},
    @{Sheet="ALC"; Cell="N55"; Value=-1224.75},
    @{Sheet="ALC"; Cell="H88"; Value=2370.25},
    @{Sheet="ALC"; Cell="I88"; Value=739.5},
    @{Sheet="ALC"; Cell="J88"; Value=2913.8333},
    @{Sheet="ALC"; Cell="K88"; Value=739.5},
    @{Sheet="ALC"; Cell="L88"; Value=2913.8333},
    @{Sheet="ALC"; Cell="M88"; Value=-333.5},
    @{Sheet="ALC"; Cell="N88"; Value=-3725.8333},
    @{Sheet="ALC"; Cell="H91"; Value=2370.25},
    @{Sheet="ALC"; Cell="I91"; Value=739.5},
    @{Sheet="ALC"; Cell="J91"; Value=2913.8333},
    @{Sheet="ALC"; Cell="K91"; Value=739.5},
    @{Sheet="ALC"; Cell="L91"; Value=2913.8333},
    @{Sheet="ALC"; Cell="M91"; Value=664.5},
    @{Sheet="ALC"; Cell="N91"; Value=-5721.8333},
    @{Sheet="ALC"; Cell="H101"; Value=609.8},
    @{Sheet="ALC"; Cell="I101"; Value=483.33334},
    @{Sheet="ALC"; Cell="J101"; Value=799.5},
    @{Sheet="ALC"; Cell="K101"; Value=1450.00002},
    @{Sheet="ALC"; Cell="L101"; Value=2398.5},
    @{Sheet="ALC"; Cell="M101"; Value=171.9999800000001},
    @{Sheet="ALC"; Cell="N101"; Value=-5642.5},
    @{Sheet="ALC"; Cell="H135"; Value=797.6},
    @{Sheet="ALC"; Cell="I135"; Value=562.15},
    @{Sheet="ALC"; Cell="J135"; Value=1739.4},
    @{Sheet="ALC"; Cell="K135"; Value=5059.349999999999},
    @{Sheet="ALC"; Cell="L135"; Value=15654.6},
    @{Sheet="ALC"; Cell="M135"; Value=-2524.349999999999},
    @{Sheet="ALC"; Cell="N135"; Value=-20724.6},
    @{Sheet="ALC"; Cell="H137"; Value=4228.0586},
    @{Sheet="ALC"; Cell="I137"; Value=1680.8334},
    @{Sheet="ALC"; Cell="J137"; Value=5617.4546},
    @{Sheet="ALC"; Cell="K137"; Value=5042.5002},
    @{Sheet="ALC"; Cell="L137"; Value=16852.3638},
    @{Sheet="ALC"; Cell="M137"; Value=-2492.5002},
    @{Sheet="ALC"; Cell="N137"; Value=-21952.3638},
    @{Sheet="ARM"; Cell="H14"; Value=4166.6665},
    @{Sheet="ARM"; Cell="I14"; Value=2500},
    @{Sheet="ARM"; Cell="K14"; Value=2500},
    @{Sheet="ARM"; Cell="M14"; Value=-2325},
    @{Sheet="ARM"; Cell="H27"; Value=4498},
    @{Sheet="ARM"; Cell="J27"; Value=4996},
    @{Sheet="ARM"; Cell="L27"; Value=4996},
    @{Sheet="ARM"; Cell="N27"; Value=-5364},
    @{Sheet="ARM"; Cell="H32"; Value=14538.657},
    @{Sheet="ARM"; Cell="I32"; Value=6187.5127},
    @{Sheet="ARM"; Cell="J32"; Value=24117.912},
    @{Sheet="ARM"; Cell="K32"; Value=6187.5127},
    @{Sheet="ARM"; Cell="L32"; Value=24117.912},
    @{Sheet="ARM"; Cell="M32"; Value=-5900.5127},
    @{Sheet="ARM"; Cell="N32"; Value=-24691.912},
    @{Sheet="ARM"; Cell="H61"; Value=1495.875},
    @{Sheet="ARM"; Cell="I61"; Value=1495.875},
    @{Sheet="ARM"; Cell="J61"; Value=0},
    @{Sheet="ARM"; Cell="K61"; Value=1495.875},
    @{Sheet="ARM"; Cell="L61"; Value=0},
    @{Sheet="ARM"; Cell="M61"; Value=-1283.875},
    @{Sheet="ARM"; Cell="H63"; Value=7434.2},
    @{Sheet="ARM"; Cell="J63"; Value=8749.75},
    @{Sheet="ARM"; Cell="L63"; Value=8749.75},
    @{Sheet="ARM"; Cell="N63"; Value=-10121.75},
    @{Sheet="ARM"; Cell="H66"; Value=7434.2},
    @{Sheet="ARM"; Cell="J66"; Value=8749.75},
    @{Sheet="ARM"; Cell="L66"; Value=43748.75},
    @{Sheet="ARM"; Cell="N66"; Value=-50612.75},
    @{Sheet="ARM"; Cell="H132"; Value=2425.3667},
    @{Sheet="ARM"; Cell="I132"; Value=1865.8846},
    @{Sheet="ARM"; Cell="K132"; Value=5597.6538},
    @{Sheet="ARM"; Cell="M132"; Value=-3067.6538},
    @{Sheet="ARM"; Cell="H136"; Value=1495.875},
    @{Sheet="ARM"; Cell="I136"; Value=1495.875},
    @{Sheet="ARM"; Cell="J136"; Value=0},
    @{Sheet="ARM"; Cell="K136"; Value=4487.625},
    @{Sheet="ARM"; Cell="L136"; Value=0},
    @{Sheet="ARM"; Cell="M136"; Value=-1937.625},
    @{Sheet="BSM"; Cell="H105"; Value=4129.4},
    @{Sheet="BSM"; Cell="I105"; Value=3591.1052},
    @{Sheet="BSM"; Cell="J105"; Value=4768.625},
    @{Sheet="BSM"; Cell="K105"; Value=3591.1052},
    @{Sheet="BSM"; Cell="L105"; Value=4768.625},
    @{Sheet="BSM"; Cell="M105"; Value=-1844.1052},
    @{Sheet="BSM"; Cell="N105"; Value=-8262.625},
    @{Sheet="BSM"; Cell="H107"; Value=2029.5834},
    @{Sheet="BSM"; Cell="I107"; Value=1636.5625},
    @{Sheet="BSM"; Cell="K107"; Value=1636.5625},
    @{Sheet="BSM"; Cell="M107"; Value=283.4375},
    @{Sheet="BSM"; Cell="H134"; Value=2478.1875},
    @{Sheet="BSM"; Cell="I134"; Value=912.5},
    @{Sheet="BSM"; Cell="K134"; Value=2737.5},
    @{Sheet="BSM"; Cell="M134"; Value=-202.5},
    @{Sheet="CRP"; Cell="H3"; Value=19996.5},
    @{Sheet="CRP"; Cell="J3"; Value=19996.5},
    @{Sheet="CRP"; Cell="L3"; Value=19996.5},
    @{Sheet="CRP"; Cell="N3"; Value=-20222.5},
    @{Sheet="CRP"; Cell="H62"; Value=54308.875},
    @{Sheet="CRP"; Cell="I62"; Value=4912},
    @{Sheet="CRP"; Cell="J62"; Value=202499.5},
    @{Sheet="CRP"; Cell="K62"; Value=4912},
    @{Sheet="CRP"; Cell="L62"; Value=202499.5},
    @{Sheet="CRP"; Cell="M62"; Value=-4288},
    @{Sheet="CRP"; Cell="N62"; Value=-203747.5},
    @{Sheet="CRP"; Cell="H65"; Value=54308.875},
    @{Sheet="CRP"; Cell="I65"; Value=4912},
    @{Sheet="CRP"; Cell="J65"; Value=202499.5},
    @{Sheet="CRP"; Cell="K65"; Value=24560},
    @{Sheet="CRP"; Cell="L65"; Value=1012497.5},
    @{Sheet="CRP"; Cell="M65"; Value=-21440},
    @{Sheet="CRP"; Cell="N65"; Value=-1018737.5},
    @{Sheet="CRP"; Cell="H129"; Value=52333.332},
    @{Sheet="CRP"; Cell="J129"; Value=52333.332},
    @{Sheet="CRP"; Cell="L129"; Value=52333.332},
    @{Sheet="CRP"; Cell="N129"; Value=-62333.332},
    @{Sheet="CRP"; Cell="H132"; Value=3524.5908},
    @{Sheet="CRP"; Cell="I132"; Value=3237.5293},
    @{Sheet="CRP"; Cell="J132"; Value=4500.6},
    @{Sheet="CRP"; Cell="K132"; Value=9712.5879},
    @{Sheet="CRP"; Cell="L132"; Value=13501.8},
    @{Sheet="CRP"; Cell="M132"; Value=-7182.5879},
    @{Sheet="CRP"; Cell="N132"; Value=-18561.8},
    @{Sheet="CRP"; Cell="H134"; Value=4077.2942},
    @{Sheet="CRP"; Cell="I134"; Value=3240.2856},
    @{Sheet="CRP"; Cell="K134"; Value=9720.856800000001},
    @{Sheet="CRP"; Cell="M134"; Value=-7185.856800000001},
    @{Sheet="CRP"; Cell="H141"; Value=84882.664},
    @{Sheet="CRP"; Cell="J141"; Value=92000},
    @{Sheet="CRP"; Cell="L141"; Value=92000},
    @{Sheet="CRP"; Cell="N141"; Value=-102360},
    @{Sheet="GSM"; Cell="H62"; Value=0},
    @{Sheet="GSM"; Cell="J62"; Value=0},
    @{Sheet="GSM"; Cell="L62"; Value=0},
    @{Sheet="GSM"; Cell="H65"; Value=0},
    @{Sheet="GSM"; Cell="J65"; Value=0},
    @{Sheet="GSM"; Cell="L65"; Value=0},
    @{Sheet="LTW"; Cell="H3"; Value=40005},
    @{Sheet="LTW"; Cell="J3"; Value=40005},
    @{Sheet="LTW"; Cell="L3"; Value=40005},
    @{Sheet="LTW"; Cell="N3"; Value=-40229},
    @{Sheet="LTW"; Cell="H15"; Value=40005},
    @{Sheet="LTW"; Cell="J15"; Value=40005},
    @{Sheet="LTW"; Cell="L15"; Value=40005},
    @{Sheet="LTW"; Cell="N15"; Value=-40345},
    @{Sheet="LTW"; Cell="H22"; Value=3200},
    @{Sheet="LTW"; Cell="J22"; Value=3200},
    @{Sheet="LTW"; Cell="L22"; Value=3200},
    @{Sheet="LTW"; Cell="N22"; Value=-3790},
    @{Sheet="LTW"; Cell="H27"; Value=3200},
    @{Sheet="LTW"; Cell="J27"; Value=3200},
    @{Sheet="LTW"; Cell="L27"; Value=3200},
    @{Sheet="LTW"; Cell="N27"; Value=-3414},
    @{Sheet="LTW"; Cell="H46"; Value=2800.7036},
    @{Sheet="LTW"; Cell="I46"; Value=1701.2222},
    @{Sheet="LTW"; Cell="J46"; Value=4999.6665},
    @{Sheet="LTW"; Cell="K46"; Value=1701.2222},
    @{Sheet="LTW"; Cell="L46"; Value=4999.6665},
    @{Sheet="LTW"; Cell="M46"; Value=-1513.2222},
    @{Sheet="LTW"; Cell="N46"; Value=-5375.6665},
    @{Sheet="LTW"; Cell="H48"; Value=925},
    @{Sheet="LTW"; Cell="I48"; Value=925},
    @{Sheet="LTW"; Cell="K48"; Value=925},
    @{Sheet="LTW"; Cell="M48"; Value=-264},
    @{Sheet="LTW"; Cell="H55"; Value=646.8946999999999},
    @{Sheet="LTW"; Cell="I55"; Value=415.93332},
    @{Sheet="LTW"; Cell="K55"; Value=415.93332},
    @{Sheet="LTW"; Cell="M55"; Value=-242.93332},
    @{Sheet="LTW"; Cell="H61"; Value=4896.231},
    @{Sheet="LTW"; Cell="J61"; Value=5748},
    @{Sheet="LTW"; Cell="L61"; Value=5748},
    @{Sheet="LTW"; Cell="N61"; Value=-6152},
    @{Sheet="LTW"; Cell="H113"; Value=4896.231},
    @{Sheet="LTW"; Cell="J113"; Value=5748},
    @{Sheet="LTW"; Cell="L113"; Value=5748},
    @{Sheet="LTW"; Cell="N113"; Value=-10088},
    @{Sheet="WVR"; Cell="H62"; Value=5472.5},
    @{Sheet="WVR"; Cell="I62"; Value=3665},
    @{Sheet="WVR"; Cell="K62"; Value=3665},
    @{Sheet="WVR"; Cell="M62"; Value=-3041},
    @{Sheet="WVR"; Cell="H65"; Value=5472.5},
    @{Sheet="WVR"; Cell="I65"; Value=3665},
    @{Sheet="WVR"; Cell="K65"; Value=18325},
    @{Sheet="WVR"; Cell="M65"; Value=-15205},
    @{Sheet="WVR"; Cell="H70"; Value=58898.5},
    @{Sheet="WVR"; Cell="I70"; Value=55595},
    @{Sheet="WVR"; Cell="J70"; Value=59999.668},
    @{Sheet="WVR"; Cell="K70"; Value=55595},
    @{Sheet="WVR"; Cell="L70"; Value=59999.668},
    @{Sheet="WVR"; Cell="M70"; Value=-55280},
    @{Sheet="WVR"; Cell="N70"; Value=-60629.668},
    @{Sheet="WVR"; Cell="H73"; Value=58898.5},
    @{Sheet="WVR"; Cell="I73"; Value=55595},
    @{Sheet="WVR"; Cell="J73"; Value=59999.668},
    @{Sheet="WVR"; Cell="K73"; Value=55595},
    @{Sheet="WVR"; Cell="L73"; Value=59999.668},
    @{Sheet="WVR"; Cell="M73"; Value=-54503},
    @{Sheet="WVR"; Cell="N73"; Value=-62183.668},
    @{Sheet="WVR"; Cell="H94"; Value=44496},
    @{Sheet="WVR"; Cell="J94"; Value=44496},
    @{Sheet="WVR"; Cell="L94"; Value=44496},
    @{Sheet="WVR"; Cell="N94"; Value=-46298}
)

$clearOps = @(
    @{Sheet="ARM"; Cell="N61"},
    @{Sheet="ARM"; Cell="N136"},
    @{Sheet="GSM"; Cell="N62"},
    @{Sheet="GSM"; Cell="N65"}
)

foreach ($chg in $setOps) {
    $ws = $wb.Worksheets.Item($chg.Sheet)
    $ws.Range($chg.Cell).Value = $chg.Value
}

foreach ($chg in $clearOps) {
    $ws = $wb.Worksheets.Item($chg.Sheet)
    $ws.Range($chg.Cell).ClearContents()
}
